$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 169; this pushes the existing rows 169-176
# down to 170-177, preserving all their data and formatting.
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row 169 with the new weekly record.
$ws.Cells.Item(169, 1).Value = 7
$ws.Cells.Item(169, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(169, 3).Value = "Ñuble"
$ws.Cells.Item(169, 4).Value = 45147
$ws.Cells.Item(169, 5).Value = 16
$ws.Cells.Item(169, 6).Value = "Fruta"
$ws.Cells.Item(169, 7).Value = 100108
$ws.Cells.Item(169, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(169, 9).Value = 100108002
$ws.Cells.Item(169, 10).Value = "Mango"
$ws.Cells.Item(169, 11).Value = "Sin especificar"
$ws.Cells.Item(169, 12).Value = "Primera"
$ws.Cells.Item(169, 13).Value = 60
$ws.Cells.Item(169, 14).Value = 8000
$ws.Cells.Item(169, 15).Value = 8000
$ws.Cells.Item(169, 16).Value = 8000
$ws.Cells.Item(169, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(169, 18).Value = "Brasil"
$ws.Cells.Item(169, 19).Value = 2000
$ws.Cells.Item(169, 20).Value = 4

# Ensure the date cell keeps the same date/time number format used by the
# other rows in column D.
$ws.Cells.Item(169, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
